$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E8").Value = "Good Morning12"
